$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New input values for the turns-ratio block (E3/E4) ---
$ws.Range("E3").Value = 311
$ws.Range("E4").Value = 16

# --- New labels in column D and computed ratios/currents in column E ---
# Insertion order matters for the shared-strings table: L2, then N3/N1, then L3.
$ws.Range("D20").Value = "L2"
$ws.Range("D9").Value = "N3/N1"
$ws.Range("D21").Value = "L3"

# D10 reuses the existing "N2/N1" label already used elsewhere on the sheet.
$ws.Range("D10").Value = "N2/N1"

# E9: N3/N1 turns ratio = E4/E3
$ws.Range("E9").Formula = "=E4/E3"
$ws.Range("E9").NumberFormat = "0.00E+00"
$ws.Range("E9").HorizontalAlignment = -4108

# E10: N2/N1 turns ratio (same formula style as B9) = B4/B3
$ws.Range("E10").Formula = "=B4/B3"

# E20: L2 = B21 * E10
$ws.Range("E20").Formula = "=B21*E10"
$ws.Range("E20").NumberFormat = "0.00E+00"
$ws.Range("E20").HorizontalAlignment = -4108

# E21: L3 = B21 * E9
$ws.Range("E21").Formula = "=B21*E9"
$ws.Range("E21").NumberFormat = "0.00E+00"
$ws.Range("E21").HorizontalAlignment = -4108

# Restore the active selection to E21, as in the saved workbook.
$ws.Range("E21").Select()
